# Update the "Responsible AI guardrails" lifecycle-stage callouts on slide 1
# to also call out security, per commit "Generative AI lifecycle stages".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 8 ("Rounded Rectangle 117"): standalone guardrails label box.
$shp1 = $s.Shapes.Item(8)
$shp1.TextFrame.TextRange.Text = "Secure & Responsible AI guardrails"

# Shape 17 ("Rounded Rectangle 126"): detailed activity description box.
$shp2 = $s.Shapes.Item(17)
$shp2.TextFrame.TextRange.Text = "Design & implement security controls & Responsible AI guardrails, with respect to hallucinations, toxicity, bias, fairness, safety, explainability, data privacy,  etc."
